# Update the "datetimeFigureOut" date field cached text (02-Mar-21 -> 14-Mar-21)
# everywhere it appears: the slide master and every slide layout.
$p = $ppt.ActivePresentation

$newDate = "14-Mar-21"

$m = $p.SlideMaster
for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $shp = $m.Shapes.Item($i)
    if ($shp.Name -like "*Date*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "*Date*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Remove the leftover "equals" (oval + math-equal) helper shapes from slide 1 —
# the play/pause logic is now handled by a function call instead, so these
# marker shapes are no longer needed.
$s = $p.Slides.Item(1)
$namesToDelete = @("Oval 1", "Equal 9", "Oval 14", "Equal 15")
foreach ($n in $namesToDelete) {
    $shp = $s.Shapes.Item($n)
    if ($shp -ne $null) {
        $shp.Delete()
    }
}
